$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column -> new value mapping (identical update applied to rows 2 and 3)
$colValues = @{
    "D"  = 0.07580000000000001
    "E"  = 0.0563
    "G"  = 0.1013615733736763
    "H"  = 0.1013615733736763
    "I"  = 0.1280887544125063
    "J"  = 0.1168426175871816
    "K"  = 41.8
    "L"  = 0.210791729702471
    "M"  = 2.96
    "N"  = 0.04974789915966386
    "O"  = 0.07081339712918661
    "P"  = 2.96
    "Q"  = 0.04974789915966386
    "R"  = 0.07081339712918661
    "U"  = 86.8
    "V"  = 1.458823529411765
    "W"  = 0.146615222728867
    "X"  = 0.1351317619518854
    "Y"  = 0.01148346077698167
    "Z"  = 0.9305490380103237
    "AA" = 0.10872778539436
    "AB" = 0.1351317619518854
    "AC" = -0.02640397655752534
    "AG" = -86.8
    "AJ" = 3.17948717948718
    "AK" = -0.3880196691998212
    "AP" = -3.364341085271318
}

foreach ($row in 2, 3) {
    foreach ($col in $colValues.Keys) {
        $ws.Range("$col$row").Value = $colValues[$col]
    }
}
